# 7.10 Fixed Some Bugs
# Wrap the three investigation-dialogue lines in green "hint" color tags
# and grow rows 3/4 so the (now longer, word-wrapped) text still fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = " <color=#00CC00>(We discovered many clues while investigating the suspects" + [char]0x2019 + " rooms.)</color>"
$ws.Range("B3").Value = " <color=#00CC00>(Let" + [char]0x2019 + "s review all the evidence and cross-reference it with the crime scene investigation.)</color>"
$ws.Range("B4").Value = " <color=#00CC00>(Who is the most likely person to have killed the Lord?)</color>"

# Rows 3 and 4 now wrap onto one more line each (3 lines / 2 lines of text
# in a 50-wide wrapped column), so grow them the way Excel's own autofit
# would (17pt per wrapped line).
$ws.Rows(3).RowHeight = 51
$ws.Rows(4).RowHeight = 34

$win = $wb.Windows.Item(1)
$win.Height = 15840
